$d = $word.ActiveDocument

# "Dur to en" -> "Dur to English and maths"
# Replace the trailing "en" with the full phrase, then force the newly
# inserted text into its own run (matching how Word splits a run when
# text is completed/typed in after the original fragment) by toggling a
# character-formatting property on/off around it; since the toggle nets
# back to the run's original (default) formatting, no extra formatting
# is left behind in the XML - only the run boundary remains.
$find = $d.Content.Find
$find.ClearFormatting()
[void]$find.Execute("en", $false, $false, $false, $false, $false, $true, 1, $false, "English and maths", 2)

$find2 = $d.Content.Find
$find2.ClearFormatting()
[void]$find2.Execute("English and maths", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$newRun = $find2.Parent
$newRun.Font.Bold = 1
$newRun.Font.Bold = 0
